$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nppa"
$ws.Range("C2").Value = "Npr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.051938
$ws.Range("H2").Value = 0.155814
$ws.Range("I2").Value = 0.1172837182974765
$ws.Range("J2").Value = 0.1172837182974765
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.536685
$ws.Range("N2").Value = 4.610055
$ws.Range("O2").Value = 0.4822880013826122
$ws.Range("P2").Value = 0.4822880013826122
$ws.Range("Q2").Value = 0.07981234553000001
$ws.Range("R2").Value = 0.71831110977
$ws.Range("S2").Value = 0.05656453009241122
$ws.Range("T2").Value = 0.05656453009241123
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nppa"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.051938
$ws.Range("H3").Value = 0.155814
$ws.Range("I3").Value = 0.1172837182974765
$ws.Range("J3").Value = 0.1172837182974765
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.195417
$ws.Range("N3").Value = 3.586251
$ws.Range("O3").Value = 0.3751811696924212
$ws.Range("P3").Value = 0.3751811696924212
$ws.Range("Q3").Value = 0.06208756814600001
$ws.Range("R3").Value = 0.558788113314
$ws.Range("S3").Value = 0.04400264261672363
$ws.Range("T3").Value = 0.04400264261672363
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nppa"
$ws.Range("C4").Value = "Npr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.051938
$ws.Range("H4").Value = 0.155814
$ws.Range("I4").Value = 0.1172837182974765
$ws.Range("J4").Value = 0.1172837182974765
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4541373333333333
$ws.Range("N4").Value = 1.362412
$ws.Range("O4").Value = 0.1425308289249667
$ws.Range("P4").Value = 0.1425308289249667
$ws.Range("Q4").Value = 0.02358698481866667
$ws.Range("R4").Value = 0.212282863368
$ws.Range("S4").Value = 0.0167165455883416
$ws.Range("T4").Value = 0.01671654558834161
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nppa"
$ws.Range("C5").Value = "Npr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.242445
$ws.Range("H5").Value = 0.7273350000000001
$ws.Range("I5").Value = 0.5474768201053503
$ws.Range("J5").Value = 0.5474768201053503
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.536685
$ws.Range("N5").Value = 4.610055
$ws.Range("O5").Value = 0.4822880013826122
$ws.Range("P5").Value = 0.4822880013826122
$ws.Range("Q5").Value = 0.3725615948250001
$ws.Range("R5").Value = 3.353054353425
$ws.Range("S5").Value = 0.2640415013719173
$ws.Range("T5").Value = 0.2640415013719173
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nppa"
$ws.Range("C6").Value = "Npr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.242445
$ws.Range("H6").Value = 0.7273350000000001
$ws.Range("I6").Value = 0.5474768201053503
$ws.Range("J6").Value = 0.5474768201053503
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.195417
$ws.Range("N6").Value = 3.586251
$ws.Range("O6").Value = 0.3751811696924212
$ws.Range("P6").Value = 0.3751811696924212
$ws.Range("Q6").Value = 0.2898228745650001
$ws.Range("R6").Value = 2.608405871085
$ws.Range("S6").Value = 0.2054029937466126
$ws.Range("T6").Value = 0.2054029937466126
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nppa"
$ws.Range("C7").Value = "Npr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.242445
$ws.Range("H7").Value = 0.7273350000000001
$ws.Range("I7").Value = 0.5474768201053503
$ws.Range("J7").Value = 0.5474768201053503
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4541373333333333
$ws.Range("N7").Value = 1.362412
$ws.Range("O7").Value = 0.1425308289249667
$ws.Range("P7").Value = 0.1425308289249667
$ws.Range("Q7").Value = 0.11010332578
$ws.Range("R7").Value = 0.99092993202
$ws.Range("S7").Value = 0.07803232498682046
$ws.Range("T7").Value = 0.07803232498682047
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Nppa"
$ws.Range("C8").Value = "Npr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1484576666666667
$ws.Range("H8").Value = 0.445373
$ws.Range("I8").Value = 0.3352394615971734
$ws.Range("J8").Value = 0.3352394615971734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.536685
$ws.Range("N8").Value = 4.610055
$ws.Range("O8").Value = 0.4822880013826122
$ws.Range("P8").Value = 0.4822880013826122
$ws.Range("Q8").Value = 0.2281326695016667
$ws.Range("R8").Value = 2.053194025515
$ws.Range("S8").Value = 0.1616819699182837
$ws.Range("T8").Value = 0.1616819699182838
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Nppa"
$ws.Range("C9").Value = "Npr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1484576666666667
$ws.Range("H9").Value = 0.445373
$ws.Range("I9").Value = 0.3352394615971734
$ws.Range("J9").Value = 0.3352394615971734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.195417
$ws.Range("N9").Value = 3.586251
$ws.Range("O9").Value = 0.3751811696924212
$ws.Range("P9").Value = 0.3751811696924212
$ws.Range("Q9").Value = 0.1774688185136667
$ws.Range("R9").Value = 1.597219366623
$ws.Range("S9").Value = 0.125775533329085
$ws.Range("T9").Value = 0.125775533329085
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Nppa"
$ws.Range("C10").Value = "Npr3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1484576666666667
$ws.Range("H10").Value = 0.445373
$ws.Range("I10").Value = 0.3352394615971734
$ws.Range("J10").Value = 0.3352394615971734
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4541373333333333
$ws.Range("N10").Value = 1.362412
$ws.Range("O10").Value = 0.1425308289249667
$ws.Range("P10").Value = 0.1425308289249667
$ws.Range("Q10").Value = 0.0674201688528889
$ws.Range("R10").Value = 0.606781519676
$ws.Range("S10").Value = 0.04778195834980468
$ws.Range("T10").Value = 0.04778195834980468
